# Trade #59 closed at 2026-02-17 12:52:41 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up figures and appends the
# newly-closed trade (#59) as row 60 on both the "All Trades" and
# "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.08            # Current Capital
$summary.Range("B4").Value = 0.07000000000000001 # Total P&L $
$summary.Range("B5").Value = 0.02               # Total P&L %
$summary.Range("B6").Value = 59                 # Total Trades
$summary.Range("B8").Value = 20                 # Losing Trades
$summary.Range("B9").Value = 42.37              # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.08              # Capital
$status.Range("D4").Value = 59                  # Trades
$status.Range("E4").Value = 0.07000000000000001 # P&L $
$status.Range("F4").Value = 0.08                # P&L %
$status.Range("G4").Value = 42.37               # Win Rate %

# ---------------------------------------------------------------------
# 3. Append trade #59 (new row 60) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$newRow = @(59, "2026-02-17", "12:52:35", "MarketMaking", "UP", 0.4, 0.24, "CLOSED", -40, -0.16, 100.08, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 60
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2 -or $col -eq 3) {
            # Date / time columns must stay plain text ("2026-02-17",
            # "12:52:35"), not auto-converted to a date/time serial.
            # Format as Text before the write, then restore the Normal
            # style so the cell carries no extra formatting afterwards.
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$i]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newRow[$i]
        }
    }
}

Write-Output "Applied trade #59 update across Summary, Strategy Status, All Trades, MarketMaking"
